$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 63, pushing the existing rows 63:135 down to 64:136.
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with a new weekly price record (same category/
# product metadata as the row that used to sit at 63, new price stats).
$ws.Range("A63").Value = 9
$ws.Range("B63").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C63").Value = "Metropolitana"
$ws.Range("D63").Value = 44539
$ws.Range("E63").Value = 13
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100101
$ws.Range("H63").Value = "Berries"
$ws.Range("I63").Value = 100101001
$ws.Range("J63").Value = "Arándano (blue)"
$ws.Range("K63").Value = "Sin especificar"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 400
$ws.Range("N63").Value = 5000
$ws.Range("O63").Value = 5000
$ws.Range("P63").Value = 5000
$ws.Range("Q63").Value = "$/bandeja 2 kilos"
$ws.Range("R63").Value = "Provincia de Curicó"
$ws.Range("S63").Value = 2500
$ws.Range("T63").Value = 2
